$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("2025-06-16 14:44:16", "Q Learning", "LineWorld", 1, 1, 0.9, 0.1, 0.1, "")
    ,("2025-06-16 14:44:23", "Q Learning", "LineWorld", 1, 2, 0.9, 0.1, 0.1, "")
    ,("2025-06-16 14:44:25", "Q Learning", "LineWorld", 1, 3, 0.9, 0.1, 0.1, "")
    ,("2025-06-16 14:44:32", "Q Learning", "LineWorld", 1, 4, 0.9, 0.1, 0.1, "")
    ,("2025-06-25 16:27:44", "Q Learning", "LineWorld", 1, 1, 0.9, 0.1, 0.1, "")
    ,("2025-06-25 16:27:50", "Q Learning", "LineWorld", 1, 2, 0.9, 0.1, 0.1, "")
    ,("2025-06-25 16:27:53", "Q Learning", "LineWorld", 1, 3, 0.9, 0.1, 0.1, "")
    ,("2025-06-25 16:29:06", "Q Learning", "GridWorld", 1, 1, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:29:11", "Q Learning", "GridWorld", 1, 2, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:29:15", "Q Learning", "GridWorld", 1, 3, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:36:43", "Q Learning", "MontyHall LV1", 0, 1, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:36:48", "Q Learning", "MontyHall LV1", 1, 2, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:36:54", "Q Learning", "MontyHall LV1", 0, 3, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:37:00", "Q Learning", "MontyHall LV1", 0, 4, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:37:11", "Q Learning", "MontyHall LV1", 1, 5, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:37:18", "Q Learning", "MontyHall LV1", 0, 6, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:37:29", "Q Learning", "MontyHall LV1", 0, 7, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 16:59:31", "Q Learning", "LineWorld", 1, 1, 0.9, 0.1, 0.1, "")
    ,("2025-06-25 16:59:34", "Q Learning", "LineWorld", 1, 2, 0.9, 0.1, 0.1, "")
    ,("2025-06-25 17:55:50", "Q Learning", "MontyHall LV1", 0, 1, 0.9, 0.1, 0.3, 10000)
    ,("2025-06-25 17:55:54", "Q Learning", "MontyHall LV1", 0, 2, 0.9, 0.1, 0.3, 10000)
)

$startRow = 45
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $startRow + $idx
    $row = $data[$idx]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Value = $row[7]
    $ws.Range("I$r").Value = $row[8]
}
